$d = $word.ActiveDocument

$pairs = @(
    @("960÷7=", "815÷9="),
    @("911÷2=", "763÷2="),
    @("532÷3=", "305÷5="),
    @("880÷3=", "122÷3="),
    @("904÷9=", "705÷2="),
    @("843÷2=", "651÷3="),
    @("679÷5=", "937÷2="),
    @("264÷3=", "396÷3="),
    @("135÷4=", "202÷3="),
    @("250÷6=", "691÷9="),
    @("899÷5=", "960÷2="),
    @("343÷6=", "294÷3="),
    @("315÷8=", "453÷5="),
    @("909÷2=", "155÷7="),
    @("544÷9=", "512÷9="),
    @("643÷9=", "669÷4="),
    @("121÷4=", "394÷9="),
    @("861÷2=", "402÷5="),
    @("638÷7=", "159÷5="),
    @("516÷5=", "113÷7="),
    @("513÷3=", "446÷7="),
    @("411÷9=", "460÷2="),
    @("349÷5=", "810÷9="),
    @("253÷9=", "354÷5="),
    @("320÷7=", "527÷7=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
